$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 needs to become the text string "1" (a shared string), not the
# number 1 - assigning a digit-only string straight to .Value gets
# auto-coerced to a numeric cell by Excel's normal "typed input"
# parsing. Going through a text formula first (="1") guarantees a
# string result, then Copy / PasteSpecial (values only) bakes that
# string down into a plain (non-formula) cell while keeping the
# existing cell style untouched.
$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
